$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1429.5385
$ws.Range("I112").Value = 577.6667
$ws.Range("J112").Value = 1880.5294
$ws.Range("K112").Value = 1733.0001
$ws.Range("L112").Value = 5641.5882
$ws.Range("M112").Value = -625.0001
$ws.Range("N112").Value = -7857.5882

$ws.Range("H116").Value = 898.5714
$ws.Range("I116").Value = 631.6667
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 631.6667
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 2810.3333
$ws.Range("N116").Value = -9384

$ws.Range("H129").Value = 13400.112
$ws.Range("I129").Value = 577.1177
$ws.Range("J129").Value = 16860.285
$ws.Range("K129").Value = 1731.3531
$ws.Range("L129").Value = 50580.855
$ws.Range("M129").Value = 3268.6469
$ws.Range("N129").Value = -60580.855

$ws.Range("H135").Value = 1356.8966
$ws.Range("I135").Value = 1179.6296
$ws.Range("J135").Value = 3750
$ws.Range("K135").Value = 10616.6664
$ws.Range("L135").Value = 33750
$ws.Range("M135").Value = -8081.6664
$ws.Range("N135").Value = -38820

$ws.Range("H137").Value = 2223.6
$ws.Range("I137").Value = 1409.3334
$ws.Range("K137").Value = 4228.0002
$ws.Range("M137").Value = -1678.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3162.7273
$ws.Range("I61").Value = 2684.2856
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2684.2856
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2472.2856
$ws.Range("N61").Value = -4424

$ws.Range("H74").Value = 1026.3158
$ws.Range("I74").Value = 1026.3158
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1026.3158
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -152.3158000000001
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1026.3158
$ws.Range("I77").Value = 1026.3158
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5131.579000000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -763.5790000000006
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 4529.7896
$ws.Range("I132").Value = 4681.2905
$ws.Range("J132").Value = 3858.8572
$ws.Range("K132").Value = 14043.8715
$ws.Range("L132").Value = 11576.5716
$ws.Range("M132").Value = -11513.8715
$ws.Range("N132").Value = -16636.5716

$ws.Range("H136").Value = 3162.7273
$ws.Range("I136").Value = 2684.2856
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 8052.8568
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -5502.8568
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 959.0625
$ws.Range("I99").Value = 1116.75
$ws.Range("J99").Value = 801.375
$ws.Range("K99").Value = 1116.75
$ws.Range("L99").Value = 801.375
$ws.Range("M99").Value = 381.25
$ws.Range("N99").Value = -3797.375

$ws.Range("H107").Value = 929.8182
$ws.Range("I107").Value = 982.0526
$ws.Range("J107").Value = 599
$ws.Range("K107").Value = 982.0526
$ws.Range("L107").Value = 599
$ws.Range("M107").Value = 937.9474
$ws.Range("N107").Value = -4439

$ws.Range("H134").Value = 17125.738
$ws.Range("I134").Value = 22380.688
$ws.Range("J134").Value = 2288.2354
$ws.Range("K134").Value = 67142.064
$ws.Range("L134").Value = 6864.706200000001
$ws.Range("M134").Value = -64607.064
$ws.Range("N134").Value = -11934.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11767363
$ws.Range("I31").Value = 2650.9092
$ws.Range("J31").Value = 33336002
$ws.Range("K31").Value = 2650.9092
$ws.Range("L31").Value = 33336002
$ws.Range("M31").Value = -2355.9092
$ws.Range("N31").Value = -33336592

$ws.Range("H34").Value = 11767363
$ws.Range("I34").Value = 2650.9092
$ws.Range("J34").Value = 33336002
$ws.Range("K34").Value = 2650.9092
$ws.Range("L34").Value = 33336002
$ws.Range("M34").Value = -2448.9092
$ws.Range("N34").Value = -33336406

$ws.Range("H58").Value = 1431.5
$ws.Range("I58").Value = 1492
$ws.Range("J58").Value = 1250
$ws.Range("K58").Value = 1492
$ws.Range("L58").Value = 1250
$ws.Range("M58").Value = -1289
$ws.Range("N58").Value = -1656

$ws.Range("H122").Value = 783.3333
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 2100
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = 350
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 2480.394
$ws.Range("I132").Value = 2303.913
$ws.Range("J132").Value = 2886.3
$ws.Range("K132").Value = 6911.739
$ws.Range("L132").Value = 8658.900000000001
$ws.Range("M132").Value = -4381.739
$ws.Range("N132").Value = -13718.9

$ws.Range("H134").Value = 1076.4517
$ws.Range("I134").Value = 1030
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 3090
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -555
$ws.Range("N134").Value = -10320

$ws.Range("H136").Value = 1431.5
$ws.Range("I136").Value = 1492
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 4476
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -1926
$ws.Range("N136").Value = -8850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 859.2308
$ws.Range("I97").Value = 842.7273
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 842.7273
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -346.7273
$ws.Range("N97").Value = -1942

$ws.Range("H126").Value = 3354.5454
$ws.Range("I126").Value = 3840
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 11520
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -9050
$ws.Range("N126").Value = -13790

$ws.Range("H132").Value = 98205.48
$ws.Range("I132").Value = 113695.39
$ws.Range("J132").Value = 5266
$ws.Range("K132").Value = 341086.17
$ws.Range("L132").Value = 15798
$ws.Range("M132").Value = -338556.17
$ws.Range("N132").Value = -20858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1295
$ws.Range("I68").Value = 1280
$ws.Range("J68").Value = 1300
$ws.Range("K68").Value = 1280
$ws.Range("L68").Value = 1300
$ws.Range("M68").Value = -531
$ws.Range("N68").Value = -2798

$ws.Range("H71").Value = 1295
$ws.Range("I71").Value = 1280
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 6400
$ws.Range("L71").Value = 6500
$ws.Range("M71").Value = -2656
$ws.Range("N71").Value = -13988

$ws.Range("H82").Value = 2198.4614
$ws.Range("I82").Value = 2265
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 2265
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -1904
$ws.Range("N82").Value = -2122

$ws.Range("H85").Value = 2198.4614
$ws.Range("I85").Value = 2265
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 2265
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -1017
$ws.Range("N85").Value = -3896

$ws.Range("H132").Value = 3443.95
$ws.Range("I132").Value = 3169.4119
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 9508.235700000001
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -6978.235700000001
$ws.Range("N132").Value = -20058.9995

$ws.Range("H136").Value = 6854.6113
$ws.Range("I136").Value = 9861.637000000001
$ws.Range("J136").Value = 2129.2856
$ws.Range("K136").Value = 29584.911
$ws.Range("L136").Value = 6387.8568
$ws.Range("M136").Value = -27034.911
$ws.Range("N136").Value = -11487.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 999.3333
$ws.Range("I96").Value = 999.2
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 999.2
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 373.8
$ws.Range("N96").Value = -3746
